$d = $word.ActiveDocument
$sr = $d.StoryRanges
$i = 0
foreach ($story in $sr) {
    $i = $i + 1
    if ($story.InlineShapes.Count -gt 0) {
        $s = $story.InlineShapes.Item(1)
        $s.Select()
        $sel = $word.Selection
        Write-Output ("Story " + $i + " Selection.InlineShapes.Count=" + $sel.InlineShapes.Count)
        if ($sel.InlineShapes.Count -gt 0) {
            $shp = $sel.InlineShapes.Item(1).ConvertToShape()
            $shp.Name = "SELTEST-" + $i
            $shp.ConvertToInlineShape() | Out-Null
        }
    }
}
